# Weekly price-list update: a new record (the latest week's price report)
# is inserted at the top of the data block (row 108) for this market /
# product combination, pushing the existing rows 108:128 down to 109:129.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 108, shifting rows 108:128 down to 109:129.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new weekly record. The other
# columns (A, B, C, E, F, G, I, R) repeat the same constant values used
# throughout this block, so copy them down from the row right below (now
# row 109, the former row 108) to keep everything consistent. Use Value2
# for reads (Value's getter mis-marshals in this host).
$ws.Range("A108").Value = $ws.Range("A109").Value2
$ws.Range("B108").Value = $ws.Range("B109").Value2
$ws.Range("C108").Value = $ws.Range("C109").Value2
$ws.Range("D108").Value = 45173
$ws.Range("D108").NumberFormat = $ws.Range("D109").NumberFormat
$ws.Range("E108").Value = $ws.Range("E109").Value2
$ws.Range("F108").Value = $ws.Range("F109").Value2
$ws.Range("G108").Value = $ws.Range("G109").Value2
$ws.Range("H108").Value = "Madrigal"
$ws.Range("I108").Value = $ws.Range("I109").Value2
$ws.Range("J108").Value = 100
$ws.Range("K108").Value = 13000
$ws.Range("L108").Value = 13000
$ws.Range("M108").Value = 13000
$ws.Range("N108").Value = "$/caja 40 unidades"
$ws.Range("O108").Value = "Provincia de Limarí"
$ws.Range("P108").Value = 325
$ws.Range("Q108").Value = 40
$ws.Range("R108").Value = $ws.Range("R109").Value2
